# Generate Report for Handback
#
# Adds a new handed-back file (cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.md) as row 4
# on the "Overview", "zh-cn" and "de-de" sheets, growing each sheet's table by
# one row and wiring up the matching hyperlinks - mirrors what the handback
# report generator does for every newly-processed file.

$wb = $excel.ActiveWorkbook

$newFile      = "cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.md"
$newFilePath  = "e2e\cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.md"
$newXlfZhCn   = "cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.090878fb2b9a44d52ffe1853bfdde597708a66f6.zh-cn.xlf"
$newXlfDeDe   = "cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.090878fb2b9a44d52ffe1853bfdde597708a66f6.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

$xlPasteValues = -4163

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3) : File Name | Path And Name | Extension |
#                              Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $newFile
$wsOverview.Range("B4").Value = $newFilePath
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = "2016-08-17 10:41:52"

$wsOverview.Range("B4").Style = "HyperLink"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc5f8fb9deadbeef1234567890abcdef12345678/e2e/cc5f8fb9-7cfc-40ce-a5aa-1e91ff7cb0e4.md",
    $null,
    $null,
    $newFilePath
) | Out-Null

# ---------------------------------------------------------------------------
# Sheets "zh-cn" (table1) and "de-de" (table2) share the same 16-column layout:
# Source File Name | File Extension | Status | Source Path | Priority |
# Content Duplicate | Correspond Handoff File | Correspond Handoff Datetime |
# Target File | Correspond Handback File | Correspond Handback DateTime |
# Reference Tokens | To be localized | Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------
function Add-LocaleRow($ws, $xlfName, $handoffDateTime, $handbackDateTime, $hyperlinkBase) {
    $lo = $ws.ListObjects.Item(1)
    $lo.ListRows.Add() | Out-Null

    $ws.Range("A4").Value = $newFile
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = $statusInSync
    $ws.Range("D4").Value = "e2e"
    $ws.Range("E4").Value = "ht"
    $ws.Range("G4").Value = $xlfName
    $ws.Range("H4").Value = $handoffDateTime
    $ws.Range("I4").Value = $newFile
    $ws.Range("J4").Value = $xlfName
    $ws.Range("K4").Value = $handbackDateTime

    # "True"/"False"/"" are stored as plain text (shared strings) in this
    # report, not as native booleans - copy the text from an existing cell
    # that already holds the right literal so the paste keeps it as text
    # instead of Excel's smart Boolean auto-detection on a typed Value.
    $ws.Range("F3").Copy()
    $ws.Range("F4").PasteSpecial($xlPasteValues)
    $ws.Range("M2").Copy()
    $ws.Range("M4").PasteSpecial($xlPasteValues)
    $ws.Range("O2").Copy()
    $ws.Range("O4").PasteSpecial($xlPasteValues)
    $ws.Range("L2").Copy()
    $ws.Range("L4").PasteSpecial($xlPasteValues)
    $ws.Range("N2").Copy()
    $ws.Range("N4").PasteSpecial($xlPasteValues)
    $ws.Range("P2").Copy()
    $ws.Range("P4").PasteSpecial($xlPasteValues)
    $excel.CutCopyMode = $false

    $ws.Range("A4").Style = "HyperLink"
    $ws.Range("I4").Style = "HyperLink"
    $ws.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $ws.Hyperlinks.Add($ws.Range("A4"), "$hyperlinkBase/$newFile", $null, $null, $newFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I4"), "$hyperlinkBase/$newFile", $null, $null, $newFile) | Out-Null
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Add-LocaleRow $wsZhCn $newXlfZhCn "2016-08-17 10:41:47" "2016-08-17 10:42:09" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/cc5f8fb9deadbeef1234567890abcdef12345678/e2e"

$wsDeDe = $wb.Worksheets.Item("de-de")
Add-LocaleRow $wsDeDe $newXlfDeDe "2016-08-17 10:41:52" "2016-08-17 10:42:16" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cc5f8fb9deadbeef1234567890abcdef12345678/e2e"

Write-Output "Handback report row added for $newFile"
